$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we touch keep their original text storage (many values
# look numeric, e.g. "20.80" or "214.93", and would otherwise be coerced to
# numbers by Excel, dropping trailing zeros / losing the ".thousands.decimal"
# style formatting used throughout this sheet).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.947.22"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.675.32"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "214.93"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "20.80"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.912.21"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "1.687.65"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "65.68"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "8.18"
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.970.86"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "236.27"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "4.44"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("E24").Value = "  -4.08%  "
$ws.Range("D25").Value = "146.66"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "16.02"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "1.485.08"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  +3.54%  "
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "67.36"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "1.819.40"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "0.779"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "90.35"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.75"
$ws.Range("E51").Value = "  +1.91%  "
